# Auto update: 2025-12-03 03:05:24
# Refresh the drone-stock (ACHR / JOBY) screen with the latest run:
#  - new as-of date
#  - the two tickers swap row position (ACHR now row 3, JOBY now row 2)
#  - refreshed metrics / model scores
#  - judgement + macro signal text updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- date column (A) -------------------------------------------------
# "2025-12-03" reads as a date literal to Excel's smart-entry, which would
# silently convert it to a date serial number. Route it through a text
# formula + Copy/PasteSpecial(values) round-trip so it lands back in the
# cell as plain text (matching the other text columns) instead of a date.
$ws.Range("A2").Formula = '="2025-12-03"'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("A3").Formula = '="2025-12-03"'
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

# --- row 2 now carries the Joby Aviation data ------------------------
$ws.Range("B2").Value = "Joby Aviation, Inc."
$ws.Range("C2").Value = "JOBY"
$ws.Range("D2").Value = 14.08
$ws.Range("E2").Value = 32.1
$ws.Range("F2").Value = 1.34
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 57.8
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 65.32892478746797
$ws.Range("O2").Value = "🟢 상승 우위 (다소 완화)"

# --- row 3 now carries the Archer Aviation data -----------------------
$ws.Range("B3").Value = "Archer Aviation Inc."
$ws.Range("C3").Value = "ACHR"
$ws.Range("D3").Value = 7.75
$ws.Range("E3").Value = 35.6
$ws.Range("F3").Value = 4.23
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 56
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 51
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 65.32892478746797
$ws.Range("O3").Value = "🟢 상승 우위 (다소 완화)"
